$p = $ppt.ActivePresentation

# --- Slide 2 ("Overview" wheel: Access/Manage/Analyze/Present) ---
$s2 = $p.Slides.Item(2)

# Rename the "Access" label to "Transfer"
$grp2 = $s2.Shapes.Item("Group 61")
$tb2 = $grp2.GroupItems.Item("TextBox 53")
$tb2.TextFrame.TextRange.Text = "Transfer"

# Add speaker-note text identifying this slide
$np2 = $s2.NotesPage
$notesBody2 = $np2.Shapes.Item(2)
$notesBody2.TextFrame.TextRange.Text = "competencies_overview"

# --- Slide 3 ("Access" detail: Database/Spreadsheet/Raw/SAS) ---
$s3 = $p.Slides.Item(3)

# Rename the "Access" label to "Transfer"
$grp3 = $s3.Shapes.Item("Group 2")
$tb3 = $grp3.GroupItems.Item("TextBox 5")
$tb3.TextFrame.TextRange.Text = "Transfer"

# Add speaker-note text identifying this slide
$np3 = $s3.NotesPage
$notesBody3 = $np3.Shapes.Item(2)
$notesBody3.TextFrame.TextRange.Text = "competencies_transfer"
